# Apply updates to "广州-漫展信息.xlsx" matching gh-pages output regenerated at 456a3b4
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 267
$wsExhibit.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202408/s1XJCd5n1724988622705.jpeg"
$wsExhibit.Range("F6").Value = 432
$wsExhibit.Range("F7").Value = 643
$wsExhibit.Range("F10").Value = 369
$wsExhibit.Range("F11").Value = 173
$wsExhibit.Range("F12").Value = 746
$wsExhibit.Range("F14").Value = 1884
$wsExhibit.Range("F15").Value = 399
$wsExhibit.Range("F16").Value = 4867
$wsExhibit.Range("F17").Value = 403
$wsExhibit.Range("F18").Value = 502
$wsExhibit.Range("F20").Value = 68
$wsExhibit.Range("F21").Value = 160

# --- Sheet "演出" (sheet2) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F13").Value = 105

# --- Sheet "本地生活" (sheet3) ---
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 5403
$wsLocal.Range("F3").Value = 346
$wsLocal.Range("F4").Value = 323

# --- Sheet "全部类型" (sheet4, aggregates all rows from the above sheets) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5403
$wsAll.Range("F4").Value = 346
$wsAll.Range("F6").Value = 323
$wsAll.Range("F7").Value = 267
$wsAll.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202408/s1XJCd5n1724988622705.jpeg"
$wsAll.Range("F18").Value = 432
$wsAll.Range("F19").Value = 643
$wsAll.Range("F23").Value = 369
$wsAll.Range("F24").Value = 173
$wsAll.Range("F27").Value = 746
$wsAll.Range("F29").Value = 105
$wsAll.Range("F30").Value = 1884
$wsAll.Range("F31").Value = 399
$wsAll.Range("F32").Value = 4867
$wsAll.Range("F34").Value = 403
$wsAll.Range("F35").Value = 502
$wsAll.Range("F37").Value = 68
$wsAll.Range("F39").Value = 160
